$p = $ppt.ActivePresentation
$s = $p.Slides.Item(23)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# [1] MAGRANI, Eduardo. A internet das coisas. Editora FGV, 2018.
# -> bold + red "MAGRANI" and "Eduardo"
$c1 = $tr.Characters(5, 7)
$c1.Font.Bold = $true
$c1.Font.Color.RGB = 255

$c2 = $tr.Characters(14, 7)
$c2.Font.Bold = $true
$c2.Font.Color.RGB = 255

# [2] Pires, Paulo F., et al. "Plataformas para a internet das coisas." ...
# -> bold + red "Pires, Paulo F., "
$c3 = $tr.Characters(70, 17)
$c3.Font.Bold = $true
$c3.Font.Color.RGB = 255

# [3] Carrion, Patrícia, and Manuela Quaresma. "Internet da Coisas (IoT): ...
# -> bold + red "Carrion, Patrícia, and Manuela Quaresma"
$c4 = $tr.Characters(234, 39)
$c4.Font.Bold = $true
$c4.Font.Color.RGB = 255

# [4] Santos, Bruno P., et al. "Internet das coisas: da teoria à prática." ...
# -> bold + red "Santos, Bruno P., "
$c5 = $tr.Characters(403, 18)
$c5.Font.Bold = $true
$c5.Font.Color.RGB = 255
